# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with the latest scraped values (GitHub Actions cron update).
#
# Note: several Price values are decimal-looking strings that must stay
# literal text (e.g. "0.120", "7.50") rather than being normalised into
# numbers (which would drop trailing zeros). Those are written with a
# leading apostrophe, exactly as a person typing them into Excel would,
# so Excel keeps them as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.618.16'
$ws.Range("E2").Value = '  -3.16%  '
$ws.Range("D3").Value = '2.617.55'
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'573.88"
$ws.Range("E5").Value = '  -3.99%  '
$ws.Range("D6").Value = "'154.94"
$ws.Range("E6").Value = '  -1.17%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.75%  '
$ws.Range("D9").Value = "'0.120"
$ws.Range("E9").Value = '  -5.48%  '
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("D11").Value = "'0.384"
$ws.Range("E11").Value = '  -3.27%  '
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("D13").Value = "'28.16"
$ws.Range("E13").Value = '  -1.89%  '
$ws.Range("D14").Value = '3.091.89'
$ws.Range("E14").Value = '  -1.26%  '
$ws.Range("D15").Value = "'0.0000183"
$ws.Range("E15").Value = '  -7.33%  '
$ws.Range("D16").Value = '63.544.51'
$ws.Range("E16").Value = '  -3.04%  '
$ws.Range("D17").Value = '2.639.08'
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("D18").Value = "'12.04"
$ws.Range("E18").Value = '  -4.36%  '
$ws.Range("D19").Value = "'4.62"
$ws.Range("E19").Value = '  -2.46%  '
$ws.Range("D20").Value = "'7.50"
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D21").Value = "'343.40"
$ws.Range("E21").Value = '  -1.59%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = "'67.09"
$ws.Range("E23").Value = '  -2.71%  '
$ws.Range("E24").Value = '  +3.30%  '
$ws.Range("E25").Value = '  -4.43%  '
$ws.Range("D26").Value = "'9.20"
$ws.Range("E26").Value = '  -4.76%  '
$ws.Range("D27").Value = "'582.95"
$ws.Range("E27").Value = '  +9.52%  '
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("E30").Value = '  -2.33%  '
$ws.Range("D31").Value = "'7.87"
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("E32").Value = '  -2.71%  '
$ws.Range("E33").Value = '  -3.82%  '
$ws.Range("D34").Value = "'6.48"
$ws.Range("E34").Value = '  +0.73%  '
$ws.Range("D35").Value = "'5.30"
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("D36").Value = "'0.408"
$ws.Range("E36").Value = '  -2.39%  '
$ws.Range("D37").Value = "'19.86"
$ws.Range("E37").Value = '  -2.54%  '
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").Value = "'153.63"
$ws.Range("E39").Value = '  -1.45%  '
$ws.Range("E40").Value = '  -3.67%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = "'41.31"
$ws.Range("E42").Value = '  -2.67%  '
$ws.Range("D43").Value = "'156.53"
$ws.Range("E43").Value = '  -2.95%  '
$ws.Range("D44").Value = "'2.37"
$ws.Range("E44").Value = '  +3.64%  '
$ws.Range("D45").Value = "'3.93"
$ws.Range("E45").Value = '  -3.08%  '
$ws.Range("D46").Value = "'0.0594"
$ws.Range("E46").Value = '  -2.03%  '
$ws.Range("D47").Value = "'22.74"
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").Value = "'0.630"
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("E49").Value = '  +2.02%  '
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("D51").Value = "'19.02"
$ws.Range("E51").Value = '  -3.74%  '
